# chore: update Sheets via scheduled runner
# Refreshes market-price derived columns (currentAveragePrice* / LevePrice* /
# LeveProfit*) on the per-job worksheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2027.5
$ws.Range("J17").Value = 2027.5
$ws.Range("L17").Value = 6082.5
$ws.Range("N17").Value = -6418.5
$ws.Range("H33").Value = 291.63635
$ws.Range("J33").Value = 236.8
$ws.Range("L33").Value = 236.8
$ws.Range("N33").Value = -694.8
$ws.Range("H40").Value = 52998.5
$ws.Range("I40").Value = 50998.5
$ws.Range("K40").Value = 50998.5
$ws.Range("M40").Value = -50823.5
$ws.Range("H70").Value = 4283.8335
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 5675.75
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 17027.25
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -17567.25
$ws.Range("H73").Value = 4283.8335
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 5675.75
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 17027.25
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -18899.25
$ws.Range("H131").Value = 3953.4
$ws.Range("I131").Value = 560.25
$ws.Range("J131").Value = 17526
$ws.Range("K131").Value = 1680.75
$ws.Range("L131").Value = 52578
$ws.Range("M131").Value = 3359.25
$ws.Range("N131").Value = -62658
$ws.Range("H132").Value = 4337.5366
$ws.Range("I132").Value = 1246.1471
$ws.Range("K132").Value = 3738.4413
$ws.Range("M132").Value = -1208.4413
$ws.Range("H138").Value = 6376.1304
$ws.Range("I138").Value = 2294.8
$ws.Range("J138").Value = 6694.9844
$ws.Range("K138").Value = 6884.400000000001
$ws.Range("L138").Value = 20084.9532
$ws.Range("M138").Value = -1744.400000000001
$ws.Range("N138").Value = -30364.9532

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6705.143
$ws.Range("I61").Value = 5895.4
$ws.Range("J61").Value = 7441.273
$ws.Range("K61").Value = 5895.4
$ws.Range("L61").Value = 7441.273
$ws.Range("M61").Value = -5683.4
$ws.Range("N61").Value = -7865.273
$ws.Range("H132").Value = 19250.771
$ws.Range("I132").Value = 29727.1
$ws.Range("J132").Value = 5282.3335
$ws.Range("K132").Value = 89181.29999999999
$ws.Range("L132").Value = 15847.0005
$ws.Range("M132").Value = -86651.29999999999
$ws.Range("N132").Value = -20907.0005
$ws.Range("H136").Value = 6705.143
$ws.Range("I136").Value = 5895.4
$ws.Range("J136").Value = 7441.273
$ws.Range("K136").Value = 17686.2
$ws.Range("L136").Value = 22323.819
$ws.Range("M136").Value = -15136.2
$ws.Range("N136").Value = -27423.819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4959.0435
$ws.Range("I20").Value = 4896.0625
$ws.Range("J20").Value = 5103
$ws.Range("K20").Value = 4896.0625
$ws.Range("L20").Value = 5103
$ws.Range("M20").Value = -4649.0625
$ws.Range("N20").Value = -5597
$ws.Range("H105").Value = 3574.8696
$ws.Range("I105").Value = 3490.9167
$ws.Range("J105").Value = 3666.4546
$ws.Range("K105").Value = 3490.9167
$ws.Range("L105").Value = 3666.4546
$ws.Range("M105").Value = -1743.9167
$ws.Range("N105").Value = -7160.4546
$ws.Range("H134").Value = 1770.7273
$ws.Range("I134").Value = 1513.2894
$ws.Range("J134").Value = 3401.1667
$ws.Range("K134").Value = 4539.8682
$ws.Range("L134").Value = 10203.5001
$ws.Range("M134").Value = -2004.8682
$ws.Range("N134").Value = -15273.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 5000
$ws.Range("J11").Value = 5000
$ws.Range("L11").Value = 5000
$ws.Range("N11").Value = -5280
$ws.Range("H31").Value = 2852.7473
$ws.Range("I31").Value = 2282.4285
$ws.Range("J31").Value = 5989.5
$ws.Range("K31").Value = 2282.4285
$ws.Range("L31").Value = 5989.5
$ws.Range("M31").Value = -1987.4285
$ws.Range("N31").Value = -6579.5
$ws.Range("H34").Value = 2852.7473
$ws.Range("I34").Value = 2282.4285
$ws.Range("J34").Value = 5989.5
$ws.Range("K34").Value = 2282.4285
$ws.Range("L34").Value = 5989.5
$ws.Range("M34").Value = -2080.4285
$ws.Range("N34").Value = -6393.5
$ws.Range("H132").Value = 17552104
$ws.Range("I132").Value = 25643498
$ws.Range("K132").Value = 76930494
$ws.Range("M132").Value = -76927964
$ws.Range("H134").Value = 2929.9333
$ws.Range("I134").Value = 1912.6471
$ws.Range("K134").Value = 5737.9413
$ws.Range("M134").Value = -3202.9413
$ws.Range("H139").Value = 67549.78
$ws.Range("I139").Value = 43699.6
$ws.Range("K139").Value = 43699.6
$ws.Range("M139").Value = -38559.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 695.625
$ws.Range("I44").Value = 142
$ws.Range("J44").Value = 1249.25
$ws.Range("K44").Value = 426
$ws.Range("L44").Value = 3747.75
$ws.Range("M44").Value = -28
$ws.Range("N44").Value = -4543.75
$ws.Range("H139").Value = 1141.6666
$ws.Range("I139").Value = 1141.6666
$ws.Range("K139").Value = 3424.9998
$ws.Range("M139").Value = 1715.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 11907259
$ws.Range("J102").Value = 6486.778
$ws.Range("L102").Value = 6486.778
$ws.Range("N102").Value = -9730.778
$ws.Range("H124").Value = 595000
$ws.Range("J124").Value = 595000
$ws.Range("L124").Value = 595000
$ws.Range("N124").Value = -604820
$ws.Range("H126").Value = 3810.7693
$ws.Range("I126").Value = 2001.8823
$ws.Range("K126").Value = 6005.6469
$ws.Range("M126").Value = -3535.6469
$ws.Range("H132").Value = 1665.75
$ws.Range("I132").Value = 867.25
$ws.Range("K132").Value = 2601.75
$ws.Range("M132").Value = -71.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 955.5294
$ws.Range("I22").Value = 823.5833
$ws.Range("J22").Value = 1272.2
$ws.Range("K22").Value = 823.5833
$ws.Range("L22").Value = 1272.2
$ws.Range("M22").Value = -528.5833
$ws.Range("N22").Value = -1862.2
$ws.Range("H27").Value = 955.5294
$ws.Range("I27").Value = 823.5833
$ws.Range("J27").Value = 1272.2
$ws.Range("K27").Value = 823.5833
$ws.Range("L27").Value = 1272.2
$ws.Range("M27").Value = -716.5833
$ws.Range("N27").Value = -1486.2
$ws.Range("H55").Value = 565.375
$ws.Range("I55").Value = 98
$ws.Range("K55").Value = 98
$ws.Range("M55").Value = 75
$ws.Range("H132").Value = 4140.3535
$ws.Range("I132").Value = 3491.3403
$ws.Range("K132").Value = 10474.0209
$ws.Range("M132").Value = -7944.0209
$ws.Range("H133").Value = 92326
$ws.Range("J133").Value = 92326
$ws.Range("L133").Value = 92326
$ws.Range("N133").Value = -97386
$ws.Range("H136").Value = 3868.5151
$ws.Range("I136").Value = 2799.6458
$ws.Range("J136").Value = 6718.8335
$ws.Range("K136").Value = 8398.937399999999
$ws.Range("L136").Value = 20156.5005
$ws.Range("M136").Value = -5848.937399999999
$ws.Range("N136").Value = -25256.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 21517.666
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 21517.666
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 21517.666
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -22071.666
$ws.Range("H132").Value = 1561.9753
$ws.Range("I132").Value = 587.8
$ws.Range("K132").Value = 1763.4
$ws.Range("M132").Value = 766.6000000000001
$ws.Range("H136").Value = 2830.3901
$ws.Range("I136").Value = 1866.1082
$ws.Range("K136").Value = 5598.3246
$ws.Range("M136").Value = -3048.3246
